# Auto-generated Excel COM-interop script
# Applies numeric corrections to per-leve market price/profit data
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H117").Value = 89999
$ws.Range("J117").Value = 89999
$ws.Range("L117").Value = 89999
$ws.Range("N117").Value = -99177
$ws.Range("H137").Value = 14142.192
$ws.Range("I137").Value = 6365.636
$ws.Range("J137").Value = 19845
$ws.Range("K137").Value = 19096.908
$ws.Range("L137").Value = 59535
$ws.Range("M137").Value = -16546.908
$ws.Range("N137").Value = -64635

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4788.2856
$ws.Range("I32").Value = 4850.327
$ws.Range("J32").Value = 1376
$ws.Range("K32").Value = 4850.327
$ws.Range("L32").Value = 1376
$ws.Range("M32").Value = -4563.327
$ws.Range("N32").Value = -1950
$ws.Range("H63").Value = 2724.3333
$ws.Range("I63").Value = 1514.6428
$ws.Range("K63").Value = 1514.6428
$ws.Range("M63").Value = -828.6428000000001
$ws.Range("H66").Value = 2724.3333
$ws.Range("I66").Value = 1514.6428
$ws.Range("K66").Value = 7573.214
$ws.Range("M66").Value = -4141.214
$ws.Range("H74").Value = 17119.5
$ws.Range("I74").Value = 18288.75
$ws.Range("K74").Value = 18288.75
$ws.Range("M74").Value = -17414.75
$ws.Range("H77").Value = 17119.5
$ws.Range("I77").Value = 18288.75
$ws.Range("K77").Value = 91443.75
$ws.Range("M77").Value = -87075.75
$ws.Range("H109").Value = 30000
$ws.Range("J109").Value = 30000
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32774
$ws.Range("H132").Value = 1132.75
$ws.Range("I132").Value = 1069.4412
$ws.Range("K132").Value = 3208.3236
$ws.Range("M132").Value = -678.3235999999997

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 30046.166
$ws.Range("I82").Value = 8952.111000000001
$ws.Range("J82").Value = 93328.336
$ws.Range("K82").Value = 8952.111000000001
$ws.Range("L82").Value = 93328.336
$ws.Range("M82").Value = -8569.111000000001
$ws.Range("N82").Value = -94094.336
$ws.Range("H85").Value = 30046.166
$ws.Range("I85").Value = 8952.111000000001
$ws.Range("J85").Value = 93328.336
$ws.Range("K85").Value = 8952.111000000001
$ws.Range("L85").Value = 93328.336
$ws.Range("M85").Value = -7626.111000000001
$ws.Range("N85").Value = -95980.336
$ws.Range("H86").Value = 240498.53
$ws.Range("I86").Value = 418214.97
$ws.Range("J86").Value = 3543.2778
$ws.Range("K86").Value = 418214.97
$ws.Range("L86").Value = 3543.2778
$ws.Range("M86").Value = -417091.97
$ws.Range("N86").Value = -5789.2778
$ws.Range("H89").Value = 240498.53
$ws.Range("I89").Value = 418214.97
$ws.Range("J89").Value = 3543.2778
$ws.Range("K89").Value = 2091074.85
$ws.Range("L89").Value = 17716.389
$ws.Range("M89").Value = -2085458.85
$ws.Range("N89").Value = -28948.389
$ws.Range("H97").Value = 11248.625
$ws.Range("I97").Value = 5074
$ws.Range("K97").Value = 5074
$ws.Range("M97").Value = -4083
$ws.Range("H99").Value = 3901
$ws.Range("J99").Value = 3300.375
$ws.Range("L99").Value = 3300.375
$ws.Range("N99").Value = -6296.375
$ws.Range("H103").Value = 29130.625
$ws.Range("J103").Value = 29130.625
$ws.Range("L103").Value = 29130.625
$ws.Range("N103").Value = -31474.625
$ws.Range("H134").Value = 11237.487
$ws.Range("I134").Value = 5698.346
$ws.Range("J134").Value = 22315.77
$ws.Range("K134").Value = 17095.038
$ws.Range("L134").Value = 66947.31
$ws.Range("M134").Value = -14560.038
$ws.Range("N134").Value = -72017.31

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2926.5833
$ws.Range("I31").Value = 1564.5
$ws.Range("J31").Value = 4288.6665
$ws.Range("K31").Value = 1564.5
$ws.Range("L31").Value = 4288.6665
$ws.Range("M31").Value = -1269.5
$ws.Range("N31").Value = -4878.6665
$ws.Range("H34").Value = 2926.5833
$ws.Range("I34").Value = 1564.5
$ws.Range("J34").Value = 4288.6665
$ws.Range("K34").Value = 1564.5
$ws.Range("L34").Value = 4288.6665
$ws.Range("M34").Value = -1362.5
$ws.Range("N34").Value = -4692.6665
$ws.Range("H132").Value = 21736.955
$ws.Range("I132").Value = 10934.637
$ws.Range("J132").Value = 43341.59
$ws.Range("K132").Value = 32803.911
$ws.Range("L132").Value = 130024.77
$ws.Range("M132").Value = -30273.911
$ws.Range("N132").Value = -135084.77
$ws.Range("H134").Value = 3363.7344
$ws.Range("I134").Value = 2535.0908
$ws.Range("J134").Value = 5186.75
$ws.Range("K134").Value = 7605.2724
$ws.Range("L134").Value = 15560.25
$ws.Range("M134").Value = -5070.2724
$ws.Range("N134").Value = -20630.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 1799.6666
$ws.Range("I130").Value = 1799.6666
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 5398.9998
$ws.Range("L130").Value = 0
$ws.Range("M130").Value = -378.9997999999996
$ws.Range("N130").ClearContents()
$ws.Range("H139").Value = 997.4815
$ws.Range("I139").Value = 997.4815
$ws.Range("K139").Value = 2992.4445
$ws.Range("M139").Value = 2147.5555

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 774.53845
$ws.Range("I97").Value = 427.3
$ws.Range("J97").Value = 1932
$ws.Range("K97").Value = 427.3
$ws.Range("L97").Value = 1932
$ws.Range("M97").Value = 68.69999999999999
$ws.Range("N97").Value = -2924
$ws.Range("H102").Value = 2138.7368
$ws.Range("I102").Value = 2138.7368
$ws.Range("K102").Value = 2138.7368
$ws.Range("M102").Value = -516.7368000000001
$ws.Range("H113").Value = 670666.3
$ws.Range("I113").Value = 670666.3
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 670666.3
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -668496.3
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 14965.947
$ws.Range("I132").Value = 14686.277
$ws.Range("J132").Value = 20000
$ws.Range("K132").Value = 44058.831
$ws.Range("L132").Value = 60000
$ws.Range("M132").Value = -41528.831
$ws.Range("N132").Value = -65060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6182.5454
$ws.Range("I40").Value = 6346.222
$ws.Range("J40").Value = 5446
$ws.Range("K40").Value = 6346.222
$ws.Range("L40").Value = 5446
$ws.Range("M40").Value = -6210.222
$ws.Range("N40").Value = -5718
$ws.Range("H61").Value = 1822.7273
$ws.Range("I61").Value = 2025.4286
$ws.Range("J61").Value = 1468
$ws.Range("K61").Value = 2025.4286
$ws.Range("L61").Value = 1468
$ws.Range("M61").Value = -1823.4286
$ws.Range("N61").Value = -1872
$ws.Range("H99").Value = 32921.777
$ws.Range("I99").Value = 33287
$ws.Range("J99").Value = 30000
$ws.Range("K99").Value = 33287
$ws.Range("L99").Value = 30000
$ws.Range("M99").Value = -30292
$ws.Range("N99").Value = -35990
$ws.Range("H113").Value = 1822.7273
$ws.Range("I113").Value = 2025.4286
$ws.Range("J113").Value = 1468
$ws.Range("K113").Value = 2025.4286
$ws.Range("L113").Value = 1468
$ws.Range("M113").Value = 144.5714
$ws.Range("N113").Value = -5808
$ws.Range("H122").Value = 4876.5713
$ws.Range("I122").Value = 3731.1667
$ws.Range("K122").Value = 11193.5001
$ws.Range("M122").Value = -8743.500100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 23441.666
$ws.Range("J54").Value = 36795
$ws.Range("L54").Value = 36795
$ws.Range("N54").Value = -37835
$ws.Range("H61").Value = 33184.11
$ws.Range("I61").Value = 33184.11
$ws.Range("K61").Value = 33184.11
$ws.Range("M61").Value = -32892.11
$ws.Range("H107").Value = 1422.3226
$ws.Range("I107").Value = 1092.9642
$ws.Range("K107").Value = 3278.8926
$ws.Range("M107").Value = -1358.8926
$ws.Range("H113").Value = 922.6
$ws.Range("I113").Value = 992.5454999999999
$ws.Range("J113").Value = 837.1111
$ws.Range("K113").Value = 2977.6365
$ws.Range("L113").Value = 2511.3333
$ws.Range("M113").Value = -807.6364999999996
$ws.Range("N113").Value = -6851.3333
$ws.Range("H132").Value = 117405.19
$ws.Range("I132").Value = 149860.08
$ws.Range("K132").Value = 449580.24
$ws.Range("M132").Value = -447050.24
$ws.Range("H136").Value = 2632879.8
$ws.Range("I136").Value = 3509542.5
$ws.Range("K136").Value = 10528627.5
$ws.Range("M136").Value = -10528627.5
